$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    'D2' = '98.657.80'
    'E2' = '  +1.54%  '
    'D3' = '3.300.75'
    'E3' = '  -0.90%  '
    'E4' = '  +0.01%  '
    'D5' = '254.48'
    'E5' = '  +1.32%  '
    'D6' = '623.45'
    'E6' = '  +0.27%  '
    'D7' = '1.48'
    'E7' = '  +32.90%  '
    'D8' = '0.400'
    'E8' = '  +4.27%  '
    'E9' = '  -0.06%  '
    'D10' = '0.957'
    'E10' = '  +21.58%  '
    'D11' = '3.300.06'
    'E11' = '  -0.83%  '
    'E12' = '  +0.54%  '
    'D13' = '39.31'
    'E13' = '  +10.97%  '
    'D14' = '98.309.25'
    'E14' = '  +1.37%  '
    'D15' = '0.0000248'
    'E15' = '  +0.69%  '
    'D16' = '3.916.81'
    'E16' = '  +0.00%  '
    'D17' = '5.47'
    'E17' = '  -0.76%  '
    'D18' = '3.298.25'
    'E18' = '  -0.76%  '
    'D19' = '3.46'
    'E19' = '  -2.72%  '
    'E20' = '  +3.70%  '
    'D21' = '6.29'
    'E21' = '  +8.40%  '
    'D22' = '484.75'
    'E22' = '  +0.46%  '
    'D23' = '9.42'
    'E23' = '  +1.86%  '
    'D24' = '0.0000202'
    'E24' = '  -3.98%  '
    'D25' = '5.60'
    'E25' = '  -1.29%  '
    'D26' = '88.64'
    'E26' = '  +1.15%  '
    'D27' = '11.99'
    'E27' = '  -0.80%  '
    'D28' = '0.303'
    'E28' = '  +25.89%  '
    'D29' = '3.476.07'
    'E29' = '  -0.57%  '
    'E30' = '  -0.07%  '
    'E31' = '  +11.96%  '
    'E32' = '  +2.26%  '
    'D33' = '10.04'
    'E33' = '  +9.24%  '
    'D34' = '0.997'
    'E34' = '  -0.24%  '
    'E35' = '  +2.67%  '
    'B36' = 'Kaspa'
    'C36' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'D36' = '0.147'
    'E36' = '  -2.00%  '
    'B37' = 'RenderToken'
    'C37' = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
    'D37' = '7.15'
    'E37' = '  -3.55%  '
    'D38' = '1.94'
    'E38' = '  -0.01%  '
    'D39' = '0.462'
    'E39' = '  +2.66%  '
    'E40' = '  +0.14%  '
    'D41' = '488.29'
    'E41' = '  -3.57%  '
    'D42' = '3.67'
    'E42' = '  +4.52%  '
    'E43' = '  -4.42%  '
    'D44' = '0.786'
    'E44' = '  -1.81%  '
    'E45' = '  -0.03%  '
    'D46' = '3.10'
    'E46' = '  -5.88%  '
    'D47' = '159.43'
    'E47' = '  -1.16%  '
    'E48' = '  +0.69%  '
    'D49' = '7.30'
    'E49' = '  +15.51%  '
    'B50' = 'Mantle'
    'C50' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D50' = '0.845'
    'E50' = '  +6.15%  '
    'B51' = 'Filecoin'
    'C51' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D51' = '4.72'
    'E51' = '  +4.50%  '
}

foreach ($key in $changes.Keys) {
    $cell = $ws.Range($key)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$key]
    $cell.Style = "Normal"
}

Write-Output "done"